$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 0.71
$ws.Range("C2").Value = 0.5091693569921298
$ws.Range("G2").Value = 0.05634472104730723
$ws.Range("H2").Value = 4.202269778071582
$ws.Range("L2").Value = 93.61

# Row 3 updates
$ws.Range("A3").Value = 493
$ws.Range("B3").Value = 22.29
$ws.Range("C3").Value = 2.879256230121726
$ws.Range("E3").Value = 22.63594512880691
$ws.Range("G3").Value = 1.776504290531661
$ws.Range("H3").Value = 66.83313373782879
$ws.Range("J3").Value = 1.392356952871433
$ws.Range("L3").Value = 169.53

# Row 4 updates
$ws.Range("A4").Value = 481
$ws.Range("B4").Value = 49.87
$ws.Range("C4").Value = 4.35991106436567
$ws.Range("E4").Value = 36.09245303066545
$ws.Range("G4").Value = 2.829153345215149
$ws.Range("H4").Value = 126.9131645197708
$ws.Range("J4").Value = 2.644024260828557
$ws.Range("L4").Value = 694.28
